# Daily update at 8 AM UTC
# Appends the next day's row of win counts to the tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 57
$prevRow = $newRow - 1

$ws.Cells.Item($newRow, 1).Value = 46006
$ws.Cells.Item($newRow, 2).Value = 122
$ws.Cells.Item($newRow, 3).Value = 137
$ws.Cells.Item($newRow, 4).Value = 128

# Match the date-formatted number format used by the other rows in column A.
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($prevRow, 1).NumberFormat
